$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7913
$ws1.Range("F5").Value = 15971
$ws1.Range("F7").Value = 594
$ws1.Range("F15").Value = 350
$ws1.Range("F20").Value = 412
$ws1.Range("F27").Value = 564

# Sheet "演出" (Performance) - sheet2
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 293

# Sheet "本地生活" (Local Life) - sheet3
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 481

# Sheet "全部类型" (All Types) - sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 481
$ws4.Range("F3").Value = 7913
$ws4.Range("F6").Value = 293
$ws4.Range("F7").Value = 15971
$ws4.Range("F9").Value = 594
$ws4.Range("F21").Value = 350
$ws4.Range("F30").Value = 412
$ws4.Range("F37").Value = 564
